# Update the timestamp embedded in the test e-mail addresses from
# 20251112_211458 to 20251112_215226, preserving the local-part
# prefixes and domain. The same e-mail text is shared (via the shared
# string table) between the "UsuariosRegistro" sheet (column C) and
# the "LoginData" sheet (column A), so both must be updated so every
# cell that displays one of these addresses reflects the new value.

$wb = $excel.ActiveWorkbook

$oldStamp = "20251112_211458"
$newStamp = "20251112_215226"

function Update-StampInRange {
    param($range)

    foreach ($cell in $range) {
        $current = $cell.Value2
        if ($current -and $current.ToString().Contains($oldStamp)) {
            $cell.Value = $current.ToString().Replace($oldStamp, $newStamp)
        }
    }
}

# Emails live in column C, rows 2-6 of UsuariosRegistro.
$wsUsuarios = $wb.Worksheets.Item("UsuariosRegistro")
Update-StampInRange $wsUsuarios.Range("C2:C6")

# The same emails are reused in column A, rows 2-3 of LoginData.
$wsLogin = $wb.Worksheets.Item("LoginData")
Update-StampInRange $wsLogin.Range("A2:A3")
